$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Round row 5 values to 2 decimal places ("custom accuracy") ---
$row5Values = @{
    2  = 15.44;  3  = 11.28;  4  = 0.72;   5  = 32.74;  6  = 26.59
    7  = 11.39;  8  = 42.97;  9  = 17.8;   10 = 7.85;   11 = 11.81
    12 = 13;     13 = 14.29;  14 = 3.74;   15 = 11.23;  16 = 16.63
    17 = 9.62;   18 = 0.43;   19 = 0.34;   20 = 169.88; 21 = 32.5
    22 = 10.95;  23 = 22.16;  24 = 11.74;  25 = 1.53;   26 = 21.08
    27 = 9.49;   28 = 8.06;   29 = 9.9;    30 = 13.78;  31 = 0.12
    32 = 38.35;  33 = 5.98;   34 = 13.28
}
foreach ($c in $row5Values.Keys) {
    $ws.Cells.Item(5, $c).Value = $row5Values[$c]
}

# --- 2. Delete row 6 entirely ---
$ws.Rows.Item(6).Delete()

# --- 3. Narrow a set of columns (raw stored width drops by 1) ---
# ColumnWidth (character units) maps to stored xlsx width as
# storedWidth = ColumnWidth + 5/6, so subtract 5/6 to target an exact
# integer stored width.
$offset = 5 / 6

$narrowTo7 = @(2,3,7,9,11,12,13,15,16,21,22,24,30,32,34)
foreach ($c in $narrowTo7) {
    $ws.Columns.Item($c).ColumnWidth = 7 - $offset
}

$narrowTo8 = @(20)
foreach ($c in $narrowTo8) {
    $ws.Columns.Item($c).ColumnWidth = 8 - $offset
}
